$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 24, shifting existing rows 24-34 down to 25-35
$ws.Rows.Item(24).EntireRow.Insert()

# Fill in the new row 24 values - most columns copy the pattern shared
# across this data block; only D, M, N, O, P, S differ per-row.
$ws.Cells.Item(24, 1).Value = 7
$ws.Cells.Item(24, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(24, 3).Value = "Ñuble"
$ws.Cells.Item(24, 4).Value = 45236
$ws.Cells.Item(24, 4).NumberFormat = $ws.Cells.Item(25, 4).NumberFormat
$ws.Cells.Item(24, 5).Value = 16
$ws.Cells.Item(24, 6).Value = "Fruta"
$ws.Cells.Item(24, 7).Value = 100107
$ws.Cells.Item(24, 8).Value = "Otros"
$ws.Cells.Item(24, 9).Value = 100107002
$ws.Cells.Item(24, 10).Value = "Chirimoya"
$ws.Cells.Item(24, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(24, 12).Value = "Primera"
$ws.Cells.Item(24, 13).Value = 100
$ws.Cells.Item(24, 14).Value = 22000
$ws.Cells.Item(24, 15).Value = 22000
$ws.Cells.Item(24, 16).Value = 22000
$ws.Cells.Item(24, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(24, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(24, 19).Value = 2200
$ws.Cells.Item(24, 20).Value = 10
